$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each coin row.
# Price cells are forced to Text format so values like "27.920.21" or
# "0.9978" are not auto-converted to numbers/dates by Excel, then the
# cell style is reset to "Normal" so no extra formatting is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.920.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.98%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.814.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.23%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9978"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.72%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "337.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.48%  "

$ws.Range("E6").Value = "  -0.58%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3925"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.83%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3488"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.63%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.11"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.90%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.202"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07591"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9953"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.75%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.17"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.99%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.522"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.62%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.813.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.198"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.33%  "

$ws.Range("E17").Value = "  +1.43%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06679"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.92%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "85.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.95%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9953"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.579"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.912.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.96%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.90%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.404"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.27%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.559"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.35%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.484"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.41%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.49%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "154.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.95%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.021.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.26%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "135.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.034"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.45%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.130"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08836"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.68%  "

$ws.Range("E35").Value = "  +0.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.555"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.32%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02434"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.610"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.53%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2232"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.41%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.265"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.578"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.63%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.84%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6580"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.99%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9947"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.870"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.50%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.161"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.05%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.62"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.18%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07206"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.82%  "

# Rows 38 and 39 swap coins (Hedera <-> TheSandbox) with updated values
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06559"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.25%  "

$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6929"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.63%  "
